$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 229:230, pushing the existing data
# (previously rows 229-328) down to rows 231-330.
$ws.Rows("229:230").Insert()

# Fill in the two new rows with the new Kiwi/Hayward price records.

# Row 229 - Calidad "Primera"
$ws.Cells.Item(229, 1).Value = 7
$ws.Cells.Item(229, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(229, 3).Value = "Ñuble"
$ws.Cells.Item(229, 4).Value = 45119
$ws.Cells.Item(229, 5).Value = 16
$ws.Cells.Item(229, 6).Value = "Fruta"
$ws.Cells.Item(229, 7).Value = 100101
$ws.Cells.Item(229, 8).Value = "Berries"
$ws.Cells.Item(229, 9).Value = 100101007
$ws.Cells.Item(229, 10).Value = "Kiwi"
$ws.Cells.Item(229, 11).Value = "Hayward"
$ws.Cells.Item(229, 12).Value = "Primera"
$ws.Cells.Item(229, 13).Value = 80
$ws.Cells.Item(229, 14).Value = 10000
$ws.Cells.Item(229, 15).Value = 10000
$ws.Cells.Item(229, 16).Value = 10000
$ws.Cells.Item(229, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(229, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(229, 19).Value = 556
$ws.Cells.Item(229, 20).Value = 18

# Row 230 - Calidad "Segunda"
$ws.Cells.Item(230, 1).Value = 7
$ws.Cells.Item(230, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(230, 3).Value = "Ñuble"
$ws.Cells.Item(230, 4).Value = 45119
$ws.Cells.Item(230, 5).Value = 16
$ws.Cells.Item(230, 6).Value = "Fruta"
$ws.Cells.Item(230, 7).Value = 100101
$ws.Cells.Item(230, 8).Value = "Berries"
$ws.Cells.Item(230, 9).Value = 100101007
$ws.Cells.Item(230, 10).Value = "Kiwi"
$ws.Cells.Item(230, 11).Value = "Hayward"
$ws.Cells.Item(230, 12).Value = "Segunda"
$ws.Cells.Item(230, 13).Value = 80
$ws.Cells.Item(230, 14).Value = 8000
$ws.Cells.Item(230, 15).Value = 8000
$ws.Cells.Item(230, 16).Value = 8000
$ws.Cells.Item(230, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(230, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(230, 19).Value = 444
$ws.Cells.Item(230, 20).Value = 18
